# Applies the "Removed is_targsim & Added mca_check" edit described in the
# commit message / OOXML diff:
#  - sheet1 (aggressive): data row updated, new "mca_check" column (I) appended
#  - sheet2 (aggressive_cont): "is_targsim" column removed, new "mca_check" column appended
#  - sheet3 (aggressive_badname): "is_targsim" column removed, new "mca_check" column appended
#  - new sheet4 "input_guidance" added at the end with a recommendation table
#
# NB: cell writes below are intentionally ordered (interleaved across sheets)
# so that newly-introduced shared strings get interned in the same sequence
# as the authoring session, keeping xl/sharedStrings.xml index-for-index
# equivalent to the target workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("aggressive")
$ws2 = $wb.Worksheets.Item("aggressive_cont")
$ws3 = $wb.Worksheets.Item("aggressive_badname")

# New sheet, appended after the last existing sheet.
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "input_guidance"

# --- 1. sheet1 A2 -> "soc.north.vpupll" -------------------------------------------------
$ws1.Range("A2").Value = "soc.north.vpupll"

# --- 2-6. sheet4 headers / row labels ---------------------------------------------------
$ws4.Range("A1").Value = "Running Environment"
$ws4.Range("A3").Value = "Targsim"
$ws4.Range("A4").Value = "Post-Si Platform"
$ws4.Range("A5").Value = "Pre-Si Platform"
$ws4.Range("B1").Value = "Recommended input"

# --- 7. sheet4 "TRUE/FALSE" cells -------------------------------------------------------
$ws4.Range("B4").Value = "TRUE/FALSE"
$ws4.Range("C4").Value = "TRUE/FALSE"
$ws4.Range("D4").Value = "TRUE/FALSE"
$ws4.Range("B5").Value = "TRUE/FALSE"
$ws4.Range("C5").Value = "TRUE/FALSE"
$ws4.Range("D5").Value = "TRUE/FALSE"

# --- 8-9. sheet1 I1/I2 -> "mca_check" / "every_failreg" --------------------------------
$ws1.Range("I1").Value = "mca_check"
$ws1.Range("I2").Value = "every_failreg"

# --- 10. sheet4 E3 (E4/E5 reuse the same string) ----------------------------------------
$ws4.Range("E3").Value = "every_10val/every_failreg"
$ws4.Range("E4").Value = "every_10val/every_failreg"
$ws4.Range("E5").Value = "every_10val/every_failreg"

# --- 11. sheet4 E2 -----------------------------------------------------------------------
$ws4.Range("E2").Value = "mca_check (it will be ignore if hang_detection is FALSE)"

# --- 12. sheet2 M2 (sheet3 G2 reuses it) --------------------------------------------------
$ws2.Range("M1").Value = "mca_check"
$ws2.Range("M2").Value = "every_10val"

# ---------------------------------------------------------------------------------------
# Remaining (non shared-string-order sensitive) cell writes
# ---------------------------------------------------------------------------------------

# sheet1 data row
$ws1.Range("E2").Value = $true
$ws1.Range("G2").Value = $true

# sheet4 remaining boolean + existing-string cells
$ws4.Range("B2").Value = "halt_detection"
$ws4.Range("C2").Value = "reset_detection"
$ws4.Range("D2").Value = "hang_detection"
$ws4.Range("B3").Value = $false
$ws4.Range("C3").Value = $false
$ws4.Range("D3").Value = $false

$ws4.Range("A1:D1").Font.Bold = $true
$ws4.Range("A2:E2").Font.Bold = $true

# sheet2: drop "is_targsim" column (shift headers left), append mca_check
$ws2.Range("I1").Value = "halt_detection"
$ws2.Range("J1").Value = "reset_detection"
$ws2.Range("K1").Value = "hang_detection"
$ws2.Range("L1").Value = "auto"

$ws2.Range("I2").Value = $false
$ws2.Range("J2").Value = $false
$ws2.Range("K2").Value = $false
$ws2.Range("L2").Value = $true

# sheet3: drop "is_targsim" column (shift headers left), append mca_check
$ws3.Range("C1").Value = "halt_detection"
$ws3.Range("D1").Value = "reset_detection"
$ws3.Range("E1").Value = "hang_detection"
$ws3.Range("F1").Value = "auto"
$ws3.Range("G1").Value = "mca_check"

$ws3.Range("C2").Value = $false
$ws3.Range("D2").Value = $false
$ws3.Range("E2").Value = $false
$ws3.Range("F2").Value = $true
$ws3.Range("G2").Value = "every_10val"

# ---------------------------------------------------------------------------------------
# Column widths for newly added columns (closest value reachable through the
# ColumnWidth->stored-width quantisation of this COM layer, which snaps to
# 1/6-character increments).
# ---------------------------------------------------------------------------------------
$ws1.Columns.Item(9).ColumnWidth = (14.109375 - 0.8333333333333334)
$ws2.Columns.Item(13).ColumnWidth = (15.88671875 - 0.8333333333333334)
$ws3.Columns.Item(7).ColumnWidth = (16.21875 - 0.8333333333333334)

$ws4.Columns.Item(1).ColumnWidth = (23.77734375 - 0.8333333333333334)
$ws4.Columns.Item(2).ColumnWidth = (19.109375 - 0.8333333333333334)
$ws4.Columns.Item(3).ColumnWidth = (16.21875 - 0.8333333333333334)
$ws4.Columns.Item(4).ColumnWidth = (16.77734375 - 0.8333333333333334)
$ws4.Columns.Item(5).ColumnWidth = (50.6640625 - 0.8333333333333334)
$ws4.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------------------
# Selections (matches target sheetView/selection in the diff)
# ---------------------------------------------------------------------------------------
[void]$ws2.Activate()
[void]$ws2.Range("M2").Select()

[void]$ws3.Activate()
[void]$ws3.Range("E9").Select()

[void]$ws4.Activate()
[void]$ws4.Range("E6").Select()

[void]$ws1.Activate()
[void]$ws1.Range("F6").Select()
